$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 93. This shifts the existing rows 93..184
# down to 94..185 (Excel preserves formatting/styles automatically on
# insert), matching the dimension growing from A1:R184 to A1:R185.
$ws.Rows.Item(93).Insert()

# Populate the newly inserted row 93 with the new weekly price observation.
# The "fixed" columns (A,B,C,E,F,G,H,I,N,O,Q,R) repeat the same values used
# throughout this subset sheet.
$ws.Range("A93").Value = 3
$ws.Range("B93").Value = "Femacal de La Calera"
$ws.Range("C93").Value = "Coquimbo"
$ws.Range("D93").Value = 44586
$ws.Range("E93").Value = 5
$ws.Range("F93").Value = 100112010
$ws.Range("G93").Value = "Achicoria"
$ws.Range("H93").Value = "Sin especificar"
$ws.Range("I93").Value = "Primera"
$ws.Range("J93").Value = 70
$ws.Range("K93").Value = 5500
$ws.Range("L93").Value = 5500
$ws.Range("M93").Value = 5500
$ws.Range("N93").Value = "$/caja 16 unidades"
$ws.Range("O93").Value = "Provincia de Quillota"
$ws.Range("P93").Value = 344
$ws.Range("Q93").Value = 16
$ws.Range("R93").Value = "Hortaliza"
